$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Update status text and timestamps to reflect "Ready for handoff" generation
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-30 15:18:30"

$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-30 15:18:26"

$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-30 15:18:30"

# Narrow the date/status columns that were widened previously.
# (Target stored width is 17.2159881591797 character-units; the COM layer
# quantizes ColumnWidth to whole pixels (1/6 character-unit steps), so
# 16.3333333333333 is the closest input that rounds to that pixel width.)
$wsOverview.Range("E:F").ColumnWidth = 16.3333333333333
$wsZhCn.Range("C:C").ColumnWidth = 16.3333333333333
$wsDeDe.Range("C:C").ColumnWidth = 16.3333333333333
